# Swap the two worksheets' tab order/names (review_info now first, hotel_info
# second) and add a new "State" column (value "Louisiana") to hotel_info,
# inserted right after "Hotel_Name" and before "City".
#
# The underlying physical sheet parts keep their position (sheet1.xml stays
# tab #1, sheet2.xml stays tab #2); what moves is the *content* + *name* of
# each tab, mirroring how Excel re-saves a sheet reorder.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # currently "hotel_info" (has the data row)
$ws2 = $wb.Worksheets.Item(2)   # currently "review_info" (header row only)

# --- Step 1: stash review_info's header row out of the way (into ws1, far
# below any real data) before we start overwriting ws2. Copy (not set)
# preserves cell types/shared-string reuse exactly. ---
$ws2.Range("A1:Y1").Copy($ws1.Range("A50"))
$ws2.Rows.Item(1).ClearContents()

# --- Step 2: move hotel_info's data (currently sitting in ws1) over to ws2,
# leaving a gap for the new "State" column between Hotel_Name (B) and
# City (originally C, now D). ---
$ws1.Range("A1:B2").Copy($ws2.Range("A1"))     # STR, Hotel_Name
$ws2.Range("C1").Value = "State"
$ws2.Range("C2").Value = "Louisiana"
$ws1.Range("C1:I2").Copy($ws2.Range("D1"))     # City ... Total_Reviews_num

# --- Step 3: clear the old hotel_info rows out of ws1. ---
$ws1.Rows.Item(1).ClearContents()
$ws1.Rows.Item(2).ClearContents()

# --- Step 4: move the stashed review_info header row back to row 1 of ws1,
# then clean up the stash. ---
$ws1.Range("A50:Y50").Copy($ws1.Range("A1"))
$ws1.Rows.Item(50).ClearContents()

# --- Step 5: rename the tabs (via a temp name to dodge the collision, since
# both final names are already taken by each other at this point). ---
$ws1.Name = "review_info_tmp"
$ws2.Name = "hotel_info"
$ws1.Name = "review_info"
